$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Per-cell updates matching the scraped coinranking.com snapshot refresh.
# Column D ("Price") holds numeric-looking text (e.g. "605.04" or "63.844.08")
# that must stay plain text, exactly as authored -- an apostrophe prefix forces
# Excel to treat the entry as text instead of auto-converting it to a number,
# and resetting the style back to Normal avoids leaving a stray text-format style
# behind on the cell (the workbook never carried one for these cells).

$ws.Range("D2").Value = "'63.844.08"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.66%  "

$ws.Range("D3").Value = "'3.320.37"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.49%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").Value = "'605.04"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.74%  "

$ws.Range("E6").Value = "  +0.40%  "

$ws.Range("E7").Value = "  +0.06%  "

$ws.Range("D8").Value = "'3.319.40"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.58%  "

$ws.Range("E9").Value = "  -0.01%  "

$ws.Range("E10").Value = "  +1.46%  "

$ws.Range("D11").Value = "'5.57"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.20%  "

$ws.Range("E12").Value = "  +0.92%  "

$ws.Range("E13").Value = "  +0.21%  "

$ws.Range("D14").Value = "'35.07"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.76%  "

$ws.Range("D15").Value = "'3.865.03"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.50%  "

$ws.Range("E16").Value = "  +0.28%  "

$ws.Range("D17").Value = "'3.318.36"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.50%  "

$ws.Range("D18").Value = "'63.923.43"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.74%  "

$ws.Range("E19").Value = "  +1.09%  "

$ws.Range("D20").Value = "'481.47"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.17%  "

$ws.Range("D21").Value = "'14.12"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.33%  "

$ws.Range("D22").Value = "'0.740"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.59%  "

$ws.Range("D23").Value = "'7.99"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.16%  "

$ws.Range("D24").Value = "'13.97"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +6.05%  "

$ws.Range("D25").Value = "'85.15"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.34%  "

$ws.Range("E26").Value = "  +0.02%  "

$ws.Range("E27").Value = "  +1.59%  "

$ws.Range("E28").Value = "  -0.07%  "

$ws.Range("D29").Value = "'8.25"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.99%  "

$ws.Range("E30").Value = "  -5.27%  "

$ws.Range("D31").Value = "'2.16"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.21%  "

$ws.Range("D32").Value = "'28.96"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.28%  "

$ws.Range("E33").Value = "  -1.41%  "

$ws.Range("E34").Value = "  -0.18%  "

$ws.Range("D35").Value = "'1.11"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.61%  "

$ws.Range("D36").Value = "'6.10"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.98%  "

$ws.Range("D37").Value = "'52.46"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.55%  "

$ws.Range("D38").Value = "'0.0₃0746"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.39%  "

$ws.Range("D39").Value = "'0.0400"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.47%  "

$ws.Range("D40").Value = "'434.96"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.97%  "

$ws.Range("D41").Value = "'3.123.94"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.71%  "

$ws.Range("E42").Value = "  +7.10%  "

$ws.Range("B43").Value = "Cosmos"
$ws.Range("C43").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D43").Value = "'8.37"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.35%  "

$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D44").Value = "'2.76"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.28%  "

$ws.Range("E45").Value = "  +0.22%  "

$ws.Range("E46").Value = "  +3.53%  "

$ws.Range("D47").Value = "'36.87"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +7.13%  "

$ws.Range("D48").Value = "'26.40"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.69%  "

$ws.Range("E50").Value = "  -2.01%  "

$ws.Range("B51").Value = "Stellar"
$ws.Range("C51").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D51").Value = "'0.114"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.59%  "
